$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows starting at row 525 - this shifts the former
# rows 525-529 down to become rows 530-534, making room for a new
# week's worth of price rows.
$ws.Range("525:529").Insert()

# Helper to write one full data row (columns A:T) given the row number.
function Set-Row {
    param($r, $a, $b, $c, $d, $e, $f, $g, $h, $i, $j, $k, $l, $m, $n, $o, $p, $q, $rr, $s, $t)
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rr
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
}

# Row 523: new data (Red Globe, Limarí)
Set-Row 523 3 "Femacal de La Calera" "Coquimbo" 44595 5 "Fruta" 100109 "Uva" 100109001 "Uva" "Red Globe" "Primera" 85 10000 10000 10000 "$/bandeja 10 kilos" "Provincia de Limarí" 1000 10

# Row 524: new data (Red Globe, San Felipe de Aconcagua)
Set-Row 524 3 "Femacal de La Calera" "Coquimbo" 44595 5 "Fruta" 100109 "Uva" 100109001 "Uva" "Red Globe" "Primera" 70 15000 15000 15000 "$/caja 15 kilos" "Provincia de San Felipe de Aconcagua" 1000 15

# Row 525: new data (Rosada pastilla, Limarí)
Set-Row 525 3 "Femacal de La Calera" "Coquimbo" 44595 5 "Fruta" 100109 "Uva" 100109001 "Uva" "Rosada pastilla" "Primera" 50 13000 13000 13000 "$/caja 15 kilos" "Provincia de Limarí" 1300 10

# Row 526: new data (Superior Seedless, San Felipe de Aconcagua)
Set-Row 526 3 "Femacal de La Calera" "Coquimbo" 44595 5 "Fruta" 100109 "Uva" 100109001 "Uva" "Superior Seedless" "Primera" 78 15000 15000 15000 "$/caja 15 kilos" "Provincia de San Felipe de Aconcagua" 1000 15

# Row 527: new data (Thompson seedless, Limarí)
Set-Row 527 3 "Femacal de La Calera" "Coquimbo" 44595 5 "Fruta" 100109 "Uva" 100109001 "Uva" "Thompson seedless" "Primera" 58 13000 13000 13000 "$/caja 15 kilos" "Provincia de Limarí" 867 15

# Row 528: carried-over data (Crimpson Seedless) - same content the old row 523 had
Set-Row 528 3 "Femacal de La Calera" "Coquimbo" 44335 5 "Fruta" 100109 "Uva" 100109001 "Uva" "Crimpson Seedless" "Primera" 65 13000 13000 13000 "$/caja 15 kilos" "Provincia de San Felipe de Aconcagua" 867 15

# Row 529: carried-over data (Red Globe) - same content the old row 524 had
Set-Row 529 3 "Femacal de La Calera" "Coquimbo" 44335 5 "Fruta" 100109 "Uva" 100109001 "Uva" "Red Globe" "Primera" 75 11000 11000 11000 "$/caja 15 kilos" "Provincia de San Felipe de Aconcagua" 733 15
